# Append three new innings rows (5,6,7) to the "Yashasvi Jaiswal" sheet,
# duplicating the existing rows 3, 4 and 2 (in that order) as per the
# source JSON re-scrape. Columns G:K hold numeric-looking values that must
# stay text (matching the original sheet's "numberStoredAsText" data), so
# those cells get their number format switched to Text before the value
# is poked in; the plain text columns (A:F) are never ambiguous and are
# set directly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "batsman" column in the source data ends with a non-breaking space
# (U+00A0), matching the existing rows 2-4 byte-for-byte.
$nbsp = [char]0x00A0
$batsman = "Yashasvi Jaiswal$nbsp"

$newRows = @(
    @{ Row = 5; Values = @(" Sharjah", " September 22 2020", "Royals won by 16 runs", "Rajasthan Royals", "Chennai Super Kings", $batsman, "6", "6", "1", "0", "100.00") },
    @{ Row = 6; Values = @(" Abu Dhabi", " October 06 2020", "Mumbai won by 57 runs", "Rajasthan Royals", "Mumbai Indians", $batsman, "0", "2", "0", "0", "0.00") },
    @{ Row = 7; Values = @(" Sharjah", " October 09 2020", "Capitals won by 46 runs", "Rajasthan Royals", "Delhi Capitals", $batsman, "34", "36", "1", "2", "94.44") }
)

# Columns G,H,I,J,K (7..11) hold numeric-looking strings that need to stay text.
$textColumns = @(7, 8, 9, 10, 11)

foreach ($entry in $newRows) {
    $r = $entry.Row
    $values = $entry.Values
    for ($c = 1; $c -le $values.Length; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        if ($textColumns -contains $c) {
            $cell.NumberFormat = "@"
        }
        $cell.Value = $values[$c - 1]
    }
}
